$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CLS")

# Insert a new column before column D, shifting existing data/formatting
# from D:K right to E:L (mirrors adding a new fiscal-year column to the
# Income Statement / Balance Sheet / Cash Flow tables).
$ws.Columns("D:D").Insert()

# Copy number formatting/styles from the (now-shifted) former column D,
# which now lives in column E, into the freshly inserted column D.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# Match the column width used by the other (bestFit) data columns.
$ws.Columns("D:D").ColumnWidth = $ws.Columns("E:E").ColumnWidth

# Populate the new column with the latest fiscal-year figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 6633200
$ws.Range("D9").Value = 6201100
$ws.Range("D10").Value = 432100
$ws.Range("D12").Value = 28800
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 61200
$ws.Range("D15").Value = 15400
$ws.Range("D17").Value = 6526900
$ws.Range("D18").Value = 106300
$ws.Range("D20").Value = -24400
$ws.Range("D21").Value = 171000
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 81900
$ws.Range("D24").Value = -17000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 98900
$ws.Range("D27").Value = 98900
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 24400
$ws.Range("D33").Value = 98900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 98900
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 422000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 1211600
$ws.Range("D44").Value = 1089900
$ws.Range("D45").Value = 100000
$ws.Range("D46").Value = 2823500
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 365300
$ws.Range("D49").Value = 492100
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 56800
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 3737700
$ws.Range("D57").Value = 1126700
$ws.Range("D58").Value = 107700
$ws.Range("D59").Value = 385900
$ws.Range("D60").Value = 1620300
$ws.Range("D61").Value = 650200
$ws.Range("D62").Value = 134900
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 2405400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -1481700
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1332300
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 98900
$ws.Range("D83").Value = 89100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 33100
$ws.Range("D91").Value = -82200
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -545600
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 419300
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -93200

Write-Output "Inserted new fiscal-year column D and populated values"
